# Apply the "Serviced by " column addition to the Card11 sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card11")

# Fix existing header text: "Correction " -> "Correction" (drop trailing space)
$ws.Range("N1").Value = "Correction"

# New header cell for the added column, formatted the same as the other
# header cells (bold, bordered, centered) by copying N1's formatting.
$ws.Range("O1").Value = "Serviced by "
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Existing N column data rows were blank placeholders; they become "nan"
# like the rest of the table's empty-value cells, and the new O column
# data rows start out blank (same as the other newly-created cells were).
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 14).Value = "nan"
    # Touch the style (without changing it) so the otherwise-empty cell is
    # still materialized in the sheet, matching the blank placeholder cells
    # used elsewhere in this column.
    $ws.Cells.Item($r, 15).Style = "Normal"
}
